$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.7967818721210449
$ws.Range("C2").Value = 0.1891467090428272
$ws.Range("D2").Value = 0.5609414031378179
$ws.Range("E2").Value = 0.1980660361244198
$ws.Range("G2").Value = 1.296034972618742
$ws.Range("H2").Value = 1.248097522473557
$ws.Range("J2").Value = 0.08498178180349036
$ws.Range("K2").Value = 0.3726252896001654
$ws.Range("L2").Value = 0.3571867946786966
$ws.Range("O2").Value = 5.193077014063135
# Row 3
$ws.Range("B3").Value = 0.7574157545826665
$ws.Range("C3").Value = 0.1893730532542826
$ws.Range("D3").Value = 0.5585613759143229
$ws.Range("E3").Value = 0.1983633079422518
$ws.Range("G3").Value = 1.303484082826756
$ws.Range("H3").Value = 1.255972145868085
$ws.Range("J3").Value = 0.08545722778968212
$ws.Range("K3").Value = 0.337933004434035
$ws.Range("L3").Value = 0.3527172189400289
$ws.Range("O3").Value = 5.225052079252421
# Row 4
$ws.Range("B4").Value = 0.7334971162608213
$ws.Range("C4").Value = 0.1895361652208152
$ws.Range("D4").Value = 0.5573358725887516
$ws.Range("E4").Value = 0.198617608128135
$ws.Range("G4").Value = 1.308689898276015
$ws.Range("H4").Value = 1.261250921833486
$ws.Range("J4").Value = 0.08576855039648112
$ws.Range("K4").Value = 0.316673838697227
$ws.Range("L4").Value = 0.3501121410265711
$ws.Range("O4").Value = 5.246943038683185
# Row 5
$ws.Range("B5").Value = 0.7238142884423553
$ws.Range("C5").Value = 0.1896087305709173
$ws.Range("D5").Value = 0.5568959204414057
$ws.Range("E5").Value = 0.1987393240166
$ws.Range("G5").Value = 1.310970292100471
$ws.Range("H5").Value = 1.263513766174071
$ws.Range("J5").Value = 0.08590030326738862
$ws.Range("K5").Value = 0.3080217003477088
$ws.Range("L5").Value = 0.3490857071604765
$ws.Range("O5").Value = 5.256431968782522
# Row 6
$ws.Range("B6").Value = 0.7222103630167851
$ws.Range("C6").Value = 0.1896211489610025
$ws.Range("D6").Value = 0.5568264616704823
$ws.Range("E6").Value = 0.1987606283544014
$ws.Range("G6").Value = 1.311358553406997
$ws.Range("H6").Value = 1.263896259764081
$ws.Range("J6").Value = 0.08592247611791581
$ws.Range("K6").Value = 0.3065857093821478
$ws.Range("L6").Value = 0.3489173960493517
$ws.Range("O6").Value = 5.258041926024859
# Row 7
$ws.Range("B7").Value = 0.7333662692531391
$ws.Range("C7").Value = 0.1895371191452888
$ws.Range("D7").Value = 0.557329698333092
$ws.Range("E7").Value = 0.198619176348565
$ws.Range("G7").Value = 1.30872000868527
$ws.Range("H7").Value = 1.261280986913199
$ws.Range("J7").Value = 0.08577030746538927
$ws.Range("K7").Value = 0.31655710704824
$ws.Range("L7").Value = 0.3500981556486025
$ws.Range("O7").Value = 5.247068708676153
# Row 8
$ws.Range("B8").Value = 0.7831565108019731
$ws.Range("C8").Value = 0.1892197590775915
$ws.Range("D8").Value = 0.5600719059984982
$ws.Range("E8").Value = 0.1981536614154358
$ws.Range("G8").Value = 1.298472277981105
$ws.Range("H8").Value = 1.250720670434859
$ws.Range("J8").Value = 0.08514169512050085
$ws.Range("K8").Value = 0.360654997324076
$ws.Range("L8").Value = 0.3556168625700025
$ws.Range("O8").Value = 5.203633552294946
# Row 9
$ws.Range("B9").Value = 0.8827688991979414
$ws.Range("C9").Value = 0.1887877840680261
$ws.Range("D9").Value = 0.5673151057351333
$ws.Range("E9").Value = 0.1978087304553533
$ws.Range("G9").Value = 1.283389708969366
$ws.Range("H9").Value = 1.233527429701908
$ws.Range("J9").Value = 0.08406249118404485
$ws.Range("K9").Value = 0.4474435605927738
$ws.Range("L9").Value = 0.3675387766425047
$ws.Range("O9").Value = 5.136360427911455
# Row 10
$ws.Range("B10").Value = 0.9571276252125642
$ws.Range("C10").Value = 0.1885849888910585
$ws.Range("D10").Value = 0.5737678337759036
$ws.Range("E10").Value = 0.1978996481676312
$ws.Range("G10").Value = 1.275363362379451
$ws.Range("H10").Value = 1.223031793448371
$ws.Range("J10").Value = 0.08336262088687185
$ws.Range("K10").Value = 0.5113762533801491
$ws.Range("L10").Value = 0.3769625151646068
$ws.Range("O10").Value = 5.097832171843407
# Row 11
$ws.Range("B11").Value = 0.991204342453841
$ws.Range("C11").Value = 0.1885173065142567
$ws.Range("D11").Value = 0.5769477178252771
$ws.Range("E11").Value = 0.1980153911153977
$ws.Range("G11").Value = 1.27237507919267
$ws.Range("H11").Value = 1.218719486093804
$ws.Range("J11").Value = 0.08306431196495012
$ws.Range("K11").Value = 0.5404934729147328
$ws.Range("L11").Value = 0.3813928216142983
$ws.Range("O11").Value = 5.082667558269833
# Row 12
$ws.Range("B12").Value = 1.004143677612802
$ws.Range("C12").Value = 0.1884951847537693
$ws.Range("D12").Value = 0.5781868879311105
$ws.Range("E12").Value = 0.1980698808750958
$ws.Range("G12").Value = 1.271338795871117
$ws.Range("H12").Value = 1.217152879387925
$ws.Range("J12").Value = 0.08295422659793772
$ws.Range("K12").Value = 0.5515237935222501
$ws.Range("L12").Value = 0.3830909627553751
$ws.Range("O12").Value = 5.077264488101264
# Row 13
$ws.Range("B13").Value = 1.001355406614664
$ws.Range("C13").Value = 0.1884997934440662
$ws.Range("D13").Value = 0.5779184553568655
$ws.Range("E13").Value = 0.1980576719768656
$ws.Range("G13").Value = 1.271557739132732
$ws.Range("H13").Value = 1.21748732582931
$ws.Range("J13").Value = 0.08297780758421602
$ws.Range("K13").Value = 0.5491480368711166
$ws.Range("L13").Value = 0.3827243290496369
$ws.Range("O13").Value = 5.078413042831073
# Row 14
$ws.Range("B14").Value = 0.992268168066829
$ws.Range("C14").Value = 0.1885154163772143
$ws.Range("D14").Value = 0.5770489641191716
$ws.Range("E14").Value = 0.1980196605756355
$ws.Range("G14").Value = 1.272287913382684
$ws.Range("H14").Value = 1.218589270794553
$ws.Range("J14").Value = 0.08305519754504154
$ws.Range("K14").Value = 0.5414008621419555
$ws.Range("L14").Value = 0.3815321190579652
$ws.Range("O14").Value = 5.082216242841071
# Row 15
$ws.Range("B15").Value = 0.9867065350489668
$ws.Range("C15").Value = 0.188525441995246
$ws.Range("D15").Value = 0.5765209316480195
$ws.Range("E15").Value = 0.1979977647308324
$ws.Range("G15").Value = 1.272747578508827
$ws.Range("H15").Value = 1.219272884499361
$ws.Range("J15").Value = 0.08310297567189195
$ws.Range("K15").Value = 0.53665603066176
$ws.Range("L15").Value = 0.3808045189765323
$ws.Range("O15").Value = 5.084590011088721
# Row 16
$ws.Range("B16").Value = 0.95490560435627
$ws.Range("C16").Value = 0.1885899043577197
$ws.Range("D16").Value = 0.5735649287666291
$ws.Range("E16").Value = 0.1978935778191087
$ws.Range("G16").Value = 1.275571992478106
$ws.Range("H16").Value = 1.223322905612619
$ws.Range("J16").Value = 0.08338251919125206
$ws.Range("K16").Value = 0.5094740044175978
$ws.Range("L16").Value = 0.3766758577823452
$ws.Range("O16").Value = 5.098870726930187
# Row 17
$ws.Range("B17").Value = 0.935460361949481
$ws.Range("C17").Value = 0.1886357250034436
$ws.Range("D17").Value = 0.5718140390738995
$ws.Range("E17").Value = 0.1978486884846937
$ws.Range("G17").Value = 1.277474463355219
$ws.Range("H17").Value = 1.225925777897814
$ws.Range("J17").Value = 0.08355914389712282
$ws.Range("K17").Value = 0.4928069598853995
$ws.Range("L17").Value = 0.3741796860921909
$ws.Range("O17").Value = 5.108236309646657
# Row 18
$ws.Range("B18").Value = 0.9242996048734256
$ws.Range("C18").Value = 0.18866439396799
$ws.Range("D18").Value = 0.5708299961960392
$ws.Range("E18").Value = 0.1978298728476169
$ws.Range("G18").Value = 1.278631114486089
$ws.Range("H18").Value = 1.227466392229374
$ws.Range("J18").Value = 0.0836626228790518
$ws.Range("K18").Value = 0.4832237457835333
$ws.Range("L18").Value = 0.3727574600870298
$ws.Range("O18").Value = 5.113845488971975
# Row 19
$ws.Range("B19").Value = 0.9205248477898635
$ws.Range("C19").Value = 0.1886744990107374
$ws.Range("D19").Value = 0.5705007746775976
$ws.Range("E19").Value = 0.197824706152673
$ws.Range("G19").Value = 1.279033454452261
$ws.Range("H19").Value = 1.227995493732763
$ws.Range("J19").Value = 0.08369798377728976
$ws.Range("K19").Value = 0.4799796075464542
$ws.Range("L19").Value = 0.3722782427420839
$ws.Range("O19").Value = 5.115782855352279
# Row 20
$ws.Range("B20").Value = 0.9375279012341764
$ws.Range("C20").Value = 0.1886306079963163
$ws.Range("D20").Value = 0.5719980425555633
$ws.Range("E20").Value = 0.1978527424534562
$ws.Range("G20").Value = 1.277265484006548
$ws.Range("H20").Value = 1.225644195157585
$ws.Range("J20").Value = 0.0835401464338883
$ws.Range("K20").Value = 0.4945808654980226
$ws.Range("L20").Value = 0.374444011036644
$ws.Range("O20").Value = 5.107216317247548
# Row 21
$ws.Range("B21").Value = 0.9949363608249087
$ws.Range("C21").Value = 0.1885107325332598
$ws.Range("D21").Value = 0.5773034055573874
$ws.Range("E21").Value = 0.1980305364109505
$ws.Range("G21").Value = 1.272070856819226
$ws.Range("H21").Value = 1.21826380232298
$ws.Range("J21").Value = 0.08303238819697434
$ws.Range("K21").Value = 0.5436762843944223
$ws.Range("L21").Value = 0.3818817453238381
$ws.Range("O21").Value = 5.081089940296039
# Row 22
$ws.Range("B22").Value = 1.032661034519037
$ws.Range("C22").Value = 0.1884528233858944
$ws.Range("D22").Value = 0.58097480680064
$ws.Range("E22").Value = 0.1982088554820791
$ws.Range("G22").Value = 1.269231427007071
$ws.Range("H22").Value = 1.213827103465064
$ws.Range("J22").Value = 0.08271730903409136
$ws.Range("K22").Value = 0.5757875338196072
$ws.Range("L22").Value = 0.3868620440449888
$ws.Range("O22").Value = 5.065993254009754
# Row 23
$ws.Range("B23").Value = 1.012508181872192
$ws.Range("C23").Value = 0.1884818688633914
$ws.Range("D23").Value = 0.5789966883294255
$ws.Range("E23").Value = 0.198108011065397
$ws.Range("G23").Value = 1.270696055663095
$ws.Range("H23").Value = 1.216159690752761
$ws.Range("J23").Value = 0.08288394075554351
$ws.Range("K23").Value = 0.5586471117510712
$ws.Range("L23").Value = 0.3841930944682872
$ws.Range("O23").Value = 5.07386968900272
# Row 24
$ws.Range("B24").Value = 0.9365931088134687
$ws.Range("C24").Value = 0.1886329141479948
$ws.Range("D24").Value = 0.5719147842787322
$ws.Range("E24").Value = 0.1978508878729457
$ws.Range("G24").Value = 1.277359767653124
$ws.Range("H24").Value = 1.225771361061192
$ws.Range("J24").Value = 0.08354872915874267
$ws.Range("K24").Value = 0.4937788861770969
$ws.Range("L24").Value = 0.3743244696735104
$ws.Range("O24").Value = 5.107676755609646
# Row 25
$ws.Range("B25").Value = 0.8556128843821966
$ws.Range("C25").Value = 0.1888844183868201
$ws.Range("D25").Value = 0.5651565348100149
$ws.Range("E25").Value = 0.1978414373652626
$ws.Range("G25").Value = 1.28693338643005
$ws.Range("H25").Value = 1.237803014370741
$ws.Range("J25").Value = 0.08433806885215489
$ws.Range("K25").Value = 0.4239337076087111
$ws.Range("L25").Value = 0.3641964404985174
$ws.Range("O25").Value = 5.152644692058715
